$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# ALERTS sheet: append a new row (row 3) with a fall-detection alert.
# ----------------------------------------------------------------------
$alerts = $wb.Worksheets.Item("ALERTS")

# Column A holds a date-like string ("2026-02-01"). Assigning it directly
# makes Excel auto-convert it to a date serial number, so format the cell
# as Text first, then restore the original (unformatted) look by copying
# the format from the row above once the text value is in place.
$alerts.Range("A3").NumberFormat = "@"
$alerts.Range("A3").Value = "2026-02-01"
$alerts.Range("A2").Copy()
$alerts.Range("A3").PasteSpecial(-4122)

$alerts.Range("B3").Value = "17:22:29"
$alerts.Range("C3").Value = "17:00"
$alerts.Range("D3").Value = "Living Room"
$alerts.Range("E3").Value = "CRITICAL"
$alerts.Range("F3").Value = "FALL_DETECTED"

# ----------------------------------------------------------------------
# mmWave sheet: append two new rows (rows 34 and 35) of presence events.
# ----------------------------------------------------------------------
$mmwave = $wb.Worksheets.Item("mmWave")

$mmwave.Range("A34").NumberFormat = "@"
$mmwave.Range("A34").Value = "2026-02-01"
$mmwave.Range("A2").Copy()
$mmwave.Range("A34").PasteSpecial(-4122)

$mmwave.Range("B34").Value = "17:22:07"
$mmwave.Range("C34").Value = "17:00"
$mmwave.Range("D34").Value = "Living Room"
$mmwave.Range("E34").Value = "PRESENCE_DETECTED"
$mmwave.Range("F34").Value = "Active"

$mmwave.Range("A35").NumberFormat = "@"
$mmwave.Range("A35").Value = "2026-02-01"
$mmwave.Range("A2").Copy()
$mmwave.Range("A35").PasteSpecial(-4122)

$mmwave.Range("B35").Value = "17:22:33"
$mmwave.Range("C35").Value = "17:00"
$mmwave.Range("D35").Value = "Living Room"
$mmwave.Range("E35").Value = "PRESENCE_DETECTED"
$mmwave.Range("F35").Value = "Active"
